$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: property/dimension-type labels
$ws.Range("B2").Value = "iaest-measure:situacion-preferente"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:nivel-estudios-agregado"

# Row 3: role labels (dim -> medida)
$ws.Range("B3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4: datatype / URI labels
$ws.Range("B4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"

# Row 5 no longer exists in the curated output - remove it entirely
$ws.Range("A5:H5").EntireRow.Delete()
